# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Numeric-looking text values (e.g. "1.003") are written with a leading
# apostrophe so Excel keeps them as text instead of auto-converting to numbers,
# matching the original inlineStr cell content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.159.76"
$ws.Range("E2").Value = "  -3.31%  "

$ws.Range("D3").Value = "1.912.19"
$ws.Range("E3").Value = "  -4.13%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -1.07%  "

$ws.Range("D5").Value = "'327.74"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("D7").Value = "'0.4673"
$ws.Range("E7").Value = "  -6.10%  "

$ws.Range("D8").Value = "'0.4006"
$ws.Range("E8").Value = "  -4.53%  "

$ws.Range("D9").Value = "'53.16"
$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").Value = "'0.08419"
$ws.Range("E10").Value = "  -5.45%  "

$ws.Range("E11").Value = "  -4.73%  "

$ws.Range("D12").Value = "'22.08"
$ws.Range("E12").Value = "  -4.10%  "

$ws.Range("D13").Value = "1.970.45"
$ws.Range("E13").Value = "  -1.09%  "

$ws.Range("D14").Value = "'7.414"
$ws.Range("E14").Value = "  -7.39%  "

$ws.Range("D15").Value = "'6.055"
$ws.Range("E15").Value = "  -5.90%  "

$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = "  -1.00%  "

$ws.Range("E17").Value = "  -3.12%  "

$ws.Range("E18").Value = "  -3.57%  "

$ws.Range("D19").Value = "'0.06614"
$ws.Range("E19").Value = "  -2.28%  "

$ws.Range("D20").Value = "'17.93"
$ws.Range("E20").Value = "  -8.39%  "

$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  -0.92%  "

$ws.Range("D22").Value = "'5.723"
$ws.Range("E22").Value = "  -4.29%  "

$ws.Range("D23").Value = "28.183.08"
$ws.Range("E23").Value = "  -3.30%  "

$ws.Range("D24").Value = "'11.19"
$ws.Range("E24").Value = "  -6.63%  "

$ws.Range("D25").Value = "'2.289"
$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").Value = "2.151.48"
$ws.Range("E26").Value = "  -3.47%  "

$ws.Range("D27").Value = "'153.28"
$ws.Range("E27").Value = "  -2.59%  "

$ws.Range("D28").Value = "'19.98"
$ws.Range("E28").Value = "  -4.24%  "

$ws.Range("D29").Value = "'5.736"
$ws.Range("E29").Value = "  -9.16%  "

$ws.Range("D30").Value = "'2.117"
$ws.Range("E30").Value = "  -6.25%  "

$ws.Range("D31").Value = "'123.19"
$ws.Range("E31").Value = "  -3.31%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.9755"
$ws.Range("E32").Value = "  -6.96%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.09648"
$ws.Range("E33").Value = "  -2.21%  "

$ws.Range("D34").Value = "'1.439"
$ws.Range("E34").Value = "  -5.74%  "

$ws.Range("D35").Value = "'3.646"
$ws.Range("E35").Value = "  -2.61%  "

$ws.Range("D36").Value = "'5.531"
$ws.Range("E36").Value = "  -5.07%  "

$ws.Range("D37").Value = "'8.796"
$ws.Range("E37").Value = "  -3.65%  "

$ws.Range("E38").Value = "  -3.81%  "

$ws.Range("E39").Value = "  -5.23%  "

$ws.Range("D40").Value = "'0.06166"
$ws.Range("E40").Value = "  -3.74%  "

$ws.Range("D41").Value = "'0.6146"
$ws.Range("E41").Value = "  -5.37%  "

$ws.Range("E42").Value = "  -4.74%  "

$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").Value = "'0.1904"
$ws.Range("E44").Value = "  -4.06%  "

$ws.Range("D45").Value = "'1.311"
$ws.Range("E45").Value = "  -4.00%  "

$ws.Range("D46").Value = "'0.5838"
$ws.Range("E46").Value = "  -6.07%  "

$ws.Range("D47").Value = "'12.75"
$ws.Range("E47").Value = "  -5.34%  "

$ws.Range("D48").Value = "'2.016"
$ws.Range("E48").Value = "  -7.86%  "

$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("D51").Value = "'110.27"
$ws.Range("E51").Value = "  -2.54%  "
